# Update the "clause" sheet with refreshed convertible-bond data and
# clear out the now-unused trailing rows, then move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set (code, name, remain, expval, lastyear) replacing rows 2-35.
$data = @(
    @("sz128138","侨银转债",4.199,135.902668163942,3.748),
    @("sh113608","威派转债",4.2,120.723920375238,3.726),
    @("sz123126","瑞丰转债",3.4,129.716782116471,4.564),
    @("sh113624","正川转债",4.049,121.468608998681,4.192),
    @("sz128118","瀛通转债",2.996,121.259900937172,3.37),
    @("sz128125","华阳转债",4.498,117.705574534807,3.447),
    @("sz123106","正丹转债",3.2,125.385322233655,4.096),
    @("sh113593","沪工转债",3.998,120.843614167109,3.419),
    @("sh113610","灵康转债",4.579,120.175044700178,3.786),
    @("sz123082","北陆转债",4.993,120.664323445494,3.803),
    @("sh113628","晨丰转债",4.15,121.389092908628,4.515),
    @("sz123093","金陵转债",2.485,128.238159342328,3.921),
    @("sh113573","纵横转债",2.698,123.655189859657,3.162),
    @("sh118004","博瑞转债",4.65,120.735694683943,4.882),
    @("sz128117","道恩转债",3.586,122.592386634904,3.37),
    @("sz123039","开润转债",2.217,119.146895662515,2.855),
    @("sz123109","昌红转债",4.596,120.705344473899,4.118),
    @("sz123050","聚飞转债",3.746,123.0471141345,3.153),
    @("sz123059","银信转债",3.912,123.050084792508,3.408),
    @("sz123136","城市转债",4.595,120.473648830648,4.929),
    @("sz128123","国光转债",3.199,115.056489657313,3.438),
    @("sz123065","宝莱转债",2.188,120.333927690038,3.545),
    @("sz128066","亚泰转债",4.61,115.52056361862,2.164),
    @("sh110070","凌钢转债",2.171,116.727072152082,3.151),
    @("sz123087","明电转债",4.189,119.994987070715,3.825),
    @("sh113532","海环转债",4.585,111.296742887031,2.121),
    @("sz123100","朗科转债",3.797,116.516218578214,3.978),
    @("sh113574","华体转债",2.089,115.179495425113,3.115),
    @("sz128072","翔鹭转债",3.016,113.486787219243,2.504),
    @("sz128097","奥佳转债",4.583,114.259143293797,3.025),
    @("sh113565","宏辉转债",2.249,116.916828604058,3.025),
    @("sz128071","合兴转债",2.999,113.249834904438,2.493),
    @("sh110058","永鼎转债",2.61,113.236825928347,2.159),
    @("sz123048","应急转债",3.79,112.130394545648,3.142)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $row++
}

# Rows 36-45 no longer have data; clear their contents (formatting stays).
$ws.Range("A36:E45").ClearContents()

# Move the visible selection to reflect the refreshed range.
$ws.Application.Goto($ws.Range("A36:E46"), $false)
$ws.Range("A36").Activate()
